# Auto-generated edit script: refresh market-price-derived columns (H:N)
# on the leve-profit sheets, per the scheduled market-data runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1434
$ws.Range("I6").Value = 1166.6666
$ws.Range("J6").Value = 1701.3334
$ws.Range("K6").Value = 3499.9998
$ws.Range("L6").Value = 5104.0002
$ws.Range("M6").Value = -3387.9998
$ws.Range("N6").ClearContents()

# Row 31
$ws.Range("H31").Value = 405.6
$ws.Range("I31").Value = 84.333336
$ws.Range("J31").Value = 887.5
$ws.Range("K31").Value = 253.000008
$ws.Range("L31").Value = 2662.5
$ws.Range("M31").Value = -23.00000800000001
$ws.Range("N31").Value = -3122.5

# Row 33
$ws.Range("H33").Value = 127.14286
$ws.Range("I33").Value = 140
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 140
$ws.Range("L33").Value = 50
$ws.Range("M33").Value = 89
$ws.Range("N33").Value = -508

# Row 58
$ws.Range("H58").Value = 1690
$ws.Range("I58").Value = 1433
$ws.Range("J58").Value = 2057.1428
$ws.Range("K58").Value = 4299
$ws.Range("L58").Value = 6171.428400000001
$ws.Range("M58").Value = -4149
$ws.Range("N58").ClearContents()

# Row 127
$ws.Range("H127").Value = 2200
$ws.Range("I127").Value = 2200
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 6600
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -1640
$ws.Range("N127").ClearContents()

# Row 129
$ws.Range("H129").Value = 1400
$ws.Range("I129").Value = 1062.3334
$ws.Range("J129").Value = 2413
$ws.Range("K129").Value = 3187.0002
$ws.Range("L129").Value = 7239
$ws.Range("M129").Value = 1812.9998
$ws.Range("N129").Value = -17239

# Row 138
$ws.Range("H138").Value = 1333.3334
$ws.Range("I138").Value = 1333.3334
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4000.0002
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 1139.9998
$ws.Range("N138").ClearContents()

# Row 141
$ws.Range("H141").Value = 2496.4
$ws.Range("I141").Value = 1160.6666
$ws.Range("J141").Value = 4500
$ws.Range("K141").Value = 3481.9998
$ws.Range("L141").Value = 13500
$ws.Range("M141").Value = 1698.0002
$ws.Range("N141").Value = -23860


$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 123.5
$ws.Range("I5").Value = 123.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 123.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -10.5
$ws.Range("N5").ClearContents()

# Row 7
$ws.Range("H7").Value = 8827894
$ws.Range("I7").Value = 9172058
$ws.Range("J7").Value = 8001900
$ws.Range("K7").Value = 9172058
$ws.Range("L7").Value = 8001900
$ws.Range("M7").Value = -9171945
$ws.Range("N7").ClearContents()

# Row 80
$ws.Range("H80").Value = 613
$ws.Range("I80").Value = 320
$ws.Range("J80").Value = 1003.6667
$ws.Range("K80").Value = 320
$ws.Range("L80").Value = 1003.6667
$ws.Range("M80").Value = 678
$ws.Range("N80").Value = -2999.6667

# Row 83
$ws.Range("H83").Value = 613
$ws.Range("I83").Value = 320
$ws.Range("J83").Value = 1003.6667
$ws.Range("K83").Value = 1600
$ws.Range("L83").Value = 5018.3335
$ws.Range("M83").Value = 3392
$ws.Range("N83").Value = -15002.3335

# Row 88
$ws.Range("H88").Value = 16632.666
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 16632.666
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 16632.666
$ws.Range("N88").Value = -17444.666

# Row 91
$ws.Range("H91").Value = 16632.666
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 16632.666
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 16632.666
$ws.Range("N91").Value = -19440.666

# Row 95
$ws.Range("H95").Value = 5989
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 5989
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 5989
$ws.Range("N95").Value = -11481

# Row 106
$ws.Range("H106").Value = 67500
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 67500
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 67500
$ws.Range("N106").Value = -70024

# Row 134
$ws.Range("H134").Value = 1139.6
$ws.Range("I134").Value = 1139.6
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3418.8
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -883.7999999999997


$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 215
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 215
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 215
$ws.Range("N12").Value = -555

# Row 28
$ws.Range("H28").Value = 17368.5
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 17368.5
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 17368.5
$ws.Range("N28").Value = -17858.5

# Row 31
$ws.Range("H31").Value = 5699.45
$ws.Range("I31").Value = 1560.1
$ws.Range("J31").Value = 9838.799999999999
$ws.Range("K31").Value = 1560.1
$ws.Range("L31").Value = 9838.799999999999
$ws.Range("M31").Value = -1265.1
$ws.Range("N31").ClearContents()

# Row 33
$ws.Range("H33").Value = 500
$ws.Range("I33").Value = 500
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 500
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -121
$ws.Range("N33").ClearContents()

# Row 34
$ws.Range("H34").Value = 5699.45
$ws.Range("I34").Value = 1560.1
$ws.Range("J34").Value = 9838.799999999999
$ws.Range("K34").Value = 1560.1
$ws.Range("L34").Value = 9838.799999999999
$ws.Range("M34").Value = -1358.1
$ws.Range("N34").ClearContents()

# Row 105
$ws.Range("H105").Value = 569.8333
$ws.Range("I105").Value = 256
$ws.Range("J105").Value = 1197.5
$ws.Range("K105").Value = 256
$ws.Range("L105").Value = 1197.5
$ws.Range("M105").Value = 1491
$ws.Range("N105").Value = -4691.5

# Row 106
$ws.Range("H106").Value = 17835.5
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 17835.5
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 17835.5
$ws.Range("N106").Value = -20359.5


$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 142
$ws.Range("I8").Value = 142
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 426
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -287

# Row 32
$ws.Range("H32").Value = 620
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 766.6667
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 2300.0001
$ws.Range("M32").Value = -917
$ws.Range("N32").ClearContents()

# Row 52
$ws.Range("H52").Value = 1265.3334
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1265.3334
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 3796.0002
$ws.Range("N52").Value = -4328.0002

# Row 56
$ws.Range("H56").Value = 3305
$ws.Range("I56").Value = 3305
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 3305
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -2775

# Row 101
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# Row 114
$ws.Range("H114").Value = 588.7778
$ws.Range("I114").Value = 277.5
$ws.Range("J114").Value = 1211.3334
$ws.Range("K114").Value = 832.5
$ws.Range("L114").Value = 3634.0002
$ws.Range("M114").Value = 2421.5
$ws.Range("N114").Value = -10142.0002

# Row 131
$ws.Range("H131").Value = 1579.0476
$ws.Range("I131").Value = 608.3333
$ws.Range("J131").Value = 2873.3333
$ws.Range("K131").Value = 1824.9999
$ws.Range("L131").Value = 8619.999899999999
$ws.Range("M131").Value = 3215.0001
$ws.Range("N131").ClearContents()

# Row 134
$ws.Range("H134").Value = 1512.4286
$ws.Range("I134").Value = 1512.4286
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4537.2858
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 532.7142000000003

# Row 137
$ws.Range("H137").Value = 5235
$ws.Range("I137").Value = 4975
$ws.Range("J137").Value = 5365
$ws.Range("K137").Value = 14925
$ws.Range("L137").Value = 16095
$ws.Range("M137").Value = -9825
$ws.Range("N137").Value = -26295

# Row 139
$ws.Range("H139").Value = 1682.2727
$ws.Range("I139").Value = 1019.375
$ws.Range("J139").Value = 3450
$ws.Range("K139").Value = 3058.125
$ws.Range("L139").Value = 10350
$ws.Range("M139").Value = 2081.875
$ws.Range("N139").Value = -20630

# Row 140
$ws.Range("H140").Value = 1913.1875
$ws.Range("I140").Value = 1707.4
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 5122.200000000001
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = 57.79999999999927
$ws.Range("N140").ClearContents()


$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 2000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1708

# Row 37
$ws.Range("H37").Value = 2000
$ws.Range("I37").Value = 2000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -1723

# Row 106
$ws.Range("H106").Value = 17980
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 17980
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 17980
$ws.Range("N106").Value = -20504


$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2045.75
$ws.Range("I22").Value = 2496.3333
$ws.Range("J22").Value = 1775.4
$ws.Range("K22").Value = 2496.3333
$ws.Range("L22").Value = 1775.4
$ws.Range("M22").Value = -2201.3333
$ws.Range("N22").Value = -2365.4

# Row 27
$ws.Range("H27").Value = 2045.75
$ws.Range("I27").Value = 2496.3333
$ws.Range("J27").Value = 1775.4
$ws.Range("K27").Value = 2496.3333
$ws.Range("L27").Value = 1775.4
$ws.Range("M27").Value = -2389.3333
$ws.Range("N27").Value = -1989.4

# Row 101
$ws.Range("H101").Value = 5681
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 5681
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 5681
$ws.Range("N101").Value = -12171

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Row 132
$ws.Range("H132").Value = 1866.6666
$ws.Range("I132").Value = 1866.6666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5599.9998
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3069.9998

# Row 136
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1950


$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 15000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 15000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 15000
$ws.Range("N63").Value = -16248

# Row 66
$ws.Range("H66").Value = 15000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 15000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 45000
$ws.Range("N66").Value = -51240

# Row 105
$ws.Range("H105").Value = 13000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 13000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 13000
$ws.Range("N105").Value = -19988

# Row 107
$ws.Range("H107").Value = 999.5
$ws.Range("I107").Value = 999.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2998.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1078.5

